$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
  2 = @(1.02,1.05245373078953,1.050498196273333,1.059044569277268,1.069366734880161,1.044112387385348,1.057476251492177,1.053252306702465,1.061775147099514,1.072069432969611,1.058977989063518)
  3 = @(1.02,1.053637334332198,1.051376228863262,1.060171072264902,1.070742972678491,1.044440393455141,1.058309206282339,1.053942628712608,1.062715011485217,1.073260438787629,1.059812126745023)
  4 = @(1.02,1.054402877520854,1.051944032692034,1.060900039672356,1.071633988893224,1.044651272194593,1.058847312792247,1.054388334168718,1.063322623792764,1.074031057717866,1.060350997427896)
  5 = @(1.02,1.054724635368148,1.052182656192069,1.061206509655586,1.072008693578355,1.044739599723896,1.059073325370276,1.054575475133902,1.063577935357643,1.074355017825687,1.06057733096971)
  6 = @(1.02,1.054778655494706,1.052222717382886,1.061257968063565,1.072071615329744,1.044754411207038,1.059111261746782,1.054606883249049,1.063620795786465,1.074409411732313,1.060615321220222)
  7 = @(1.02,1.054407177163229,1.051947221509306,1.06090413469004,1.071638995235565,1.044652453711348,1.058850333596934,1.054390835672411,1.063326035782275,1.074035386518174,1.060354022522473)
  8 = @(1.02,1.052853804082321,1.050795002690708,1.059425267832802,1.069831739008546,1.044223521081298,1.057757933462783,1.053485807443961,1.062092891840231,1.072471948656386,1.059260071054831)
  9 = @(1.02,1.050113969821969,1.048761997676939,1.056819594534166,1.066650842085205,1.043457227138738,1.055826264655547,1.051883498579565,1.059915719071788,1.069716571201676,1.057325659056323)
  10 = @(1.02,1.048285550140683,1.047404837029391,1.055082576239466,1.064532606785848,1.042939302729622,1.054533895291822,1.05081017630435,1.05846135240616,1.067879266235344,1.056031454379866)
  11 = @(1.02,1.047493354695147,1.046816728360977,1.054330430323906,1.063615916473064,1.042713353882109,1.053973179682862,1.050344190305014,1.057830884098362,1.067083574329227,1.055469942490425)
  12 = @(1.02,1.047199023801694,1.046598210131615,1.054051047116832,1.063275492134962,1.042629172735515,1.053764736472306,1.050170916378651,1.057596590477187,1.066787997482015,1.055261203266601)
  13 = @(1.02,1.047262162176364,1.046645086148847,1.054110975916907,1.063348510950394,1.04264724136755,1.053809455896681,1.050208092658125,1.05764685224306,1.066851400728212,1.055305986197689)
  14 = @(1.02,1.0474690267403,1.046798666978861,1.054307336487411,1.063587775372001,1.042706400621163,1.053955953147607,1.05032987122865,1.057811519560941,1.067059142306232,1.055452691491514)
  15 = @(1.02,1.047596472915044,1.046893284092051,1.054428320297825,1.063735203967766,1.042742816954481,1.054046192569486,1.0504048783474,1.057912961951841,1.06718713583508,1.055543059063723)
  16 = @(1.02,1.048338115268981,1.047443858316119,1.055132493408695,1.064593455210585,1.042954262671012,1.054571084599022,1.050841076241675,1.058503179235222,1.067932070811341,1.056068696500146)
  17 = @(1.02,1.04880319829997,1.047789098023455,1.055574199781685,1.065131950666907,1.043086445621107,1.054900036937882,1.051114361456099,1.058873213491299,1.068399313215788,1.056398115989059)
  18 = @(1.02,1.049074427442299,1.047990427234811,1.055831839229621,1.065446096272931,1.04316338324669,1.05509180216513,1.051273645461022,1.059088979059296,1.068671835739063,1.056590153544925)
  19 = @(1.02,1.049166901856676,1.048059067984361,1.055919687611959,1.065553220425954,1.043189589453945,1.055157171031313,1.051327937101806,1.059162537865915,1.068764756903347,1.056655615242396)
  20 = @(1.02,1.048753304059604,1.047752061552738,1.05552680892357,1.065074170016051,1.043072280449434,1.054864754554459,1.051085052818088,1.058833519482864,1.068349183792849,1.056362783500604)
  21 = @(1.02,1.047408112351227,1.046753443150924,1.054249513284008,1.063517315908299,1.042688986703699,1.053912818021399,1.050294015617876,1.057763032180234,1.066997968200658,1.055409495108479)
  22 = @(1.019999999999999,1.046561906220734,1.046125175395483,1.053446410113629,1.062538893980068,1.042446526647155,1.05331332227867,1.049795582559542,1.057089339224466,1.066148280059126,1.054809148013033)
  23 = @(1.02,1.04701053775426,1.046458269962526,1.053872152517017,1.063057533750206,1.042575198686522,1.053631219289676,1.050059913779449,1.057446537395383,1.066598728498346,1.055127496474258)
  24 = @(1.02,1.048775849260403,1.047768796873001,1.055548222813173,1.065100278448853,1.043078681582528,1.054880697475322,1.051098296496256,1.058851455709854,1.068371835157934,1.056378749062244)
  25 = @(1.02,1.050822602260819,1.049287896485708,1.057493199760286,1.067472753617938,1.043656575035095,1.056326451108448,1.052298631850313,1.060479078605774,1.070428961840196,1.05782655583136)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $colLetters.Length; $i++) {
    $ws.Range($colLetters[$i] + $row).Value = $vals[$i]
  }
}
